$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.156.36"
$ws.Range("E2").Value = "  +5.23%  "
$ws.Range("D3").Value = "'3.543.56"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'591.15"
$ws.Range("E5").Value = "  +5.31%  "
$ws.Range("D6").Value = "'192.70"
$ws.Range("E6").Value = "  +9.19%  "
$ws.Range("E7").Value = "  +1.11%  "
$ws.Range("D8").Value = "'3.533.25"
$ws.Range("E8").Value = "  +4.25%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").Value = "'0.659"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "'59.03"
$ws.Range("E12").Value = "  +9.99%  "
$ws.Range("D13").Value = "'0.0000294"
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("D14").Value = "'9.62"
$ws.Range("E14").Value = "  +4.17%  "
$ws.Range("D15").Value = "'4.090.66"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("D16").Value = "'19.14"
$ws.Range("E16").Value = "  +4.31%  "
$ws.Range("D17").Value = "'3.542.59"
$ws.Range("E17").Value = "  +3.31%  "
$ws.Range("D18").Value = "'69.131.30"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("D19").Value = "'12.36"
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +3.17%  "
$ws.Range("D22").Value = "'494.89"
$ws.Range("E22").Value = "  +1.69%  "
$ws.Range("E23").Value = "  +13.33%  "
$ws.Range("D24").Value = "'17.31"
$ws.Range("E24").Value = "  +21.02%  "
$ws.Range("D25").Value = "'4.46"
$ws.Range("E25").Value = "  +7.97%  "
$ws.Range("D26").Value = "'91.05"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").Value = "'3.03"
$ws.Range("E27").Value = "  +4.16%  "
$ws.Range("D28").Value = "'11.15"
$ws.Range("E28").Value = "  +4.13%  "
$ws.Range("D29").Value = "'9.21"
$ws.Range("E29").Value = "  +5.42%  "
$ws.Range("D30").Value = "'31.86"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'7.46"
$ws.Range("E31").Value = "  +13.64%  "
$ws.Range("D32").Value = "'614.18"
$ws.Range("E32").Value = "  +6.54%  "
$ws.Range("D33").Value = "'12.04"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("D34").Value = "'65.22"
$ws.Range("E34").Value = "  +4.44%  "
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("D36").Value = "'0.149"
$ws.Range("E36").Value = "  +5.29%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'37.69"
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "'0.0₃0793"
$ws.Range("E39").Value = "  +7.23%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'0.394"
$ws.Range("E40").Value = "  +5.10%  "
$ws.Range("D41").Value = "'3.55"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "'3.264.78"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("D43").Value = "'2.98"
$ws.Range("E43").Value = "  +6.56%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0440"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.65"
$ws.Range("E45").Value = "  +8.40%  "
$ws.Range("D46").Value = "'3.36"
$ws.Range("E46").Value = "  +5.96%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.137"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "  +18.82%  "
$ws.Range("D49").Value = "'9.07"
$ws.Range("E49").Value = "  +7.06%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'141.14"
$ws.Range("E51").Value = "  +0.59%  "
